$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window / view sizing (workbook.xml bookViews) ---
$excel.Width = 18350
$excel.Height = 17600

# --- Sheet view: scrolled/selected position changed ---
$ws.Application.ActiveWindow.ScrollColumn = 10
$ws.Range("N14").Select()

# --- New row of data (row 8): NYC Borough boundaries dataset ---
$ws.Range("A7:O7").Copy()
$ws.Range("A8:O8").PasteSpecial(-4122)
$ws.Range("A8").Select()
$excel.CutCopyMode = 0

$ws.Range("A8").Value = "nyc_boroughs"
$ws.Range("B8").Value = "nybbwi_25d"
$ws.Range("C8").Value = "data/raw/nyc_boroughs/"
$ws.Range("D8").Value = "shp"
$ws.Range("E8").Value = "NYC Gov"
$ws.Range("G8").Value = "Annually"
$ws.Range("K8").Value = "NYC Borough Boundaries"
$ws.Range("L8").Value = 46014
$ws.Range("M8").Value = 2025
$ws.Range("N8").Value = "https://www.nyc.gov/content/planning/pages/resources/datasets/borough-boundaries"
